$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2..172 (header in row 1). We need to
# append 4 new rows (173..176) describing the "ersilia-os/ersilia" repo,
# continuing the "id" sequence in column A (172, 173, 174, 175).

$lastRow = 172
$newRowsCount = 4

# Duplicate the last data row's formatting/shape for each new row by
# copying row 172 and inserting the copy right below the growing table.
# Repeating the copy from the same source row each time means every
# inserted row starts out as an identical clone of row 172 (values +
# number formats + styles), which we then overwrite with the real data.
for ($i = 0; $i -lt $newRowsCount; $i++) {
    $targetRow = $lastRow + 1 + $i
    $ws.Rows.Item($lastRow).Copy()
    $ws.Rows.Item($targetRow).Insert(-4121)

    # The row-insert sometimes drops the thin border that decorates the
    # "id" column (column A) for data rows; restore it so the new rows
    # keep the same bold + bordered + centered look as the rest of the id
    # column.
    $ws.Cells.Item($targetRow, 1).Borders.LineStyle = 1
}

$repoUrl = "https://github.com/ersilia-os/ersilia"
$repoName = "ersilia"
$repoAuthor = "ersilia-os"
$startDate = "07/04/2020"
$values = @("0", "0", "1", "1", "1", "0", "0", "0", "0", "1")

for ($i = 0; $i -lt $newRowsCount; $i++) {
    $targetRow = $lastRow + 1 + $i
    $id = $lastRow + $i

    $ws.Cells.Item($targetRow, 1).Value = $id
    $ws.Cells.Item($targetRow, 2).Value = $repoUrl
    $ws.Cells.Item($targetRow, 3).Value = $repoName
    $ws.Cells.Item($targetRow, 4).Value = $repoAuthor

    # Force the starting-date column to stay plain text (e.g. "07/04/2020")
    # instead of being auto-parsed into a date serial number.
    $ws.Cells.Item($targetRow, 5).NumberFormat = "@"
    $ws.Cells.Item($targetRow, 5).Value = $startDate

    # Columns F..O hold "0"/"1" flags that must stay text (like the rest of
    # the column), not get auto-coerced into numbers.
    for ($col = 0; $col -lt $values.Length; $col++) {
        $cell = $ws.Cells.Item($targetRow, 6 + $col)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$col]
    }
}
